$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 816.2655584213176
$ws.Range("C2").Value = 1165.712211379992
$ws.Range("D2").Value = 1378.339868674436
$ws.Range("E2").Value = 1540.366991430833
$ws.Range("B3").Value = 857.5566018388613
$ws.Range("C3").Value = 1212.959743826953
$ws.Range("D3").Value = 1413.480483055562
$ws.Range("E3").Value = 1565.302392880995
$ws.Range("B4").Value = 741.6630058926968
$ws.Range("C4").Value = 1073.324544809446
$ws.Range("D4").Value = 1286.876200652406
$ws.Range("E4").Value = 1457.830666557993
$ws.Range("B5").Value = 917.083831534653
$ws.Range("C5").Value = 1282.391924194888
$ws.Range("D5").Value = 1456.000842794155
$ws.Range("E5").Value = 1601.191988968261
$ws.Range("B6").Value = 901.9208893662466
$ws.Range("C6").Value = 1264.196820163211
$ws.Range("D6").Value = 1438.88443609094
$ws.Range("E6").Value = 1587.64886275616
$ws.Range("B7").Value = 939.6353147012608
$ws.Range("C7").Value = 1306.110061062109
$ws.Range("D7").Value = 1490.093466445672
$ws.Range("E7").Value = 1614.422536682491
$ws.Range("B8").Value = 833.2110855707815
$ws.Range("C8").Value = 1187.12235824837
$ws.Range("D8").Value = 1396.417172492213
$ws.Range("E8").Value = 1524.716381781736
$ws.Range("B9").Value = 946.7832509100929
$ws.Range("C9").Value = 1318.465130650536
$ws.Range("D9").Value = 1503.626982662522
$ws.Range("E9").Value = 1628.197875621464
$ws.Range("B10").Value = 962.4164304111293
$ws.Range("C10").Value = 1318.515521661503
$ws.Range("D10").Value = 1487.419970879135
$ws.Range("E10").Value = 1596.381014948465
$ws.Range("B11").Value = 953.6963681437817
$ws.Range("C11").Value = 1305.448640131307
$ws.Range("D11").Value = 1472.247238057132
$ws.Range("E11").Value = 1580.998156961288
$ws.Range("B12").Value = 877.6914250533716
$ws.Range("C12").Value = 1162.266216930732
$ws.Range("D12").Value = 1268.238832989708
$ws.Range("E12").Value = 1344.132768786582
$ws.Range("B13").Value = 957.8519584561764
$ws.Range("C13").Value = 1310.593135206027
$ws.Range("D13").Value = 1476.087934504065
$ws.Range("E13").Value = 1583.751461116511
